$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 to make room for the newest data point,
# shifting all existing rows (2-31) down to (3-32).
$ws.Rows(2).Insert()

# Populate the new row 2 with the latest Apple Mobility Data entry.
$ws.Range("A2").Value = 43936
$ws.Range("B2").Value = 1647
$ws.Range("C2").Value = 130216

# The row insert copied formatting down from the header row (row 1, bold
# centered). Restore the plain/date formatting used by the rest of the
# table by copying formats from the row right below (the old row 2, now
# row 3) back onto the new row 2.
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)
